# "code remis au propre" - journal de travail: fill in end-of-day time,
# add two new journal entries (maps 5-12 ; code cleanup), and extend the
# blank template row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Duplicate the formatting of the (until now) blank row 55 down onto the
# two new rows (56 = new entry, 57 = new blank template row), so they pick
# up the same cell styles (date/time/text formats, borders, fills, ...).
$ws.Range("A55:J55").Copy()
$ws.Range("A56:J57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Row 54: end time ("Heure fin") was left empty, fill it in now ---
$ws.Range("D54").Value = 0.63194444444444442

# --- Row 55: was the empty trailing template row, now a real entry ---
$ws.Range("A55").Value = 43928              # 07/04/2020
$ws.Range("B55").Value = 7                  # Semaine
$ws.Range("C55").Value = 0.625              # 15:00
$ws.Range("D55").Value = 0.65277777777777779 # 15:40
$ws.Range("G55").Value = "MA-20"
$ws.Range("H55").Value = "Maps"
$ws.Range("I55").Value = "nouvelles maps (5 -12)"

# --- Row 56: brand-new entry ---
$ws.Range("A56").Value = 43928              # 07/04/2020
$ws.Range("B56").Value = 7                  # Semaine
$ws.Range("C56").Value = 0.65625            # 15:45
$ws.Range("D56").Value = 0.70833333333333337 # 17:00
$ws.Range("G56").Value = "MA-20"
$ws.Range("H56").Value = "Commentaires"
$ws.Range("I56").Value = "Nettoyage du code (commentaires et fonctions)"

# Formula column ("Temps total") follows the same pattern as the rest of
# the sheet and needs to be (re)filled for the two rows that now have
# times, plus the new trailing blank template row.
$ws.Range("F55").Formula = '=IF(AND(C55<>"",D55<>""),D55-C55-E55,"")'
$ws.Range("F56").Formula = '=IF(AND(C56<>"",D56<>""),D56-C56-E56,"")'
$ws.Range("F57").Formula = '=IF(AND(C57<>"",D57<>""),D57-C57-E57,"")'

# Move/scroll the frozen view down a bit and leave the selection on the
# last edited cell, like the author did.
$excel.ActiveWindow.ScrollRow = 39
$ws.Range("I56").Select()

$wb.Save()
